# Scheduled runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 330.2
$ws.Range("I2").Value = 197.28572
$ws.Range("J2").Value = 640.3333
$ws.Range("K2").Value = 197.28572
$ws.Range("L2").Value = 640.3333
$ws.Range("M2").Value = -84.28572
$ws.Range("N2").Value = -866.3333
$ws.Range("H15").Value = 1165.1
$ws.Range("I15").Value = 1165.1
$ws.Range("K15").Value = 3495.3
$ws.Range("M15").Value = -3326.3
$ws.Range("H33").Value = 114.42857
$ws.Range("I33").Value = 114.42857
$ws.Range("K33").Value = 114.42857
$ws.Range("M33").Value = 114.57143
$ws.Range("H39").Value = 182.5625
$ws.Range("I39").Value = 78.53846
$ws.Range("K39").Value = 235.61538
$ws.Range("M39").Value = 60.38461999999998
$ws.Range("H40").Value = 4960.3887
$ws.Range("I40").Value = 3525.5833
$ws.Range("J40").Value = 7830
$ws.Range("K40").Value = 3525.5833
$ws.Range("L40").Value = 7830
$ws.Range("M40").Value = -3350.5833
$ws.Range("N40").Value = -8180
$ws.Range("H55").Value = 70.59999999999999
$ws.Range("I55").Value = 88
$ws.Range("J55").Value = 66.25
$ws.Range("K55").Value = 88
$ws.Range("L55").Value = 66.25
$ws.Range("M55").Value = 126
$ws.Range("N55").Value = -494.25
$ws.Range("H58").Value = 77
$ws.Range("I58").Value = 77
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 231
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -81
$ws.Range("N58").ClearContents()
$ws.Range("H88").Value = 870.375
$ws.Range("I88").Value = 1950
$ws.Range("J88").Value = 510.5
$ws.Range("K88").Value = 1950
$ws.Range("L88").Value = 510.5
$ws.Range("M88").Value = -1544
$ws.Range("N88").Value = -1322.5
$ws.Range("H91").Value = 870.375
$ws.Range("I91").Value = 1950
$ws.Range("J91").Value = 510.5
$ws.Range("K91").Value = 1950
$ws.Range("L91").Value = 510.5
$ws.Range("M91").Value = -546
$ws.Range("N91").Value = -3318.5
$ws.Range("H96").Value = 675.25
$ws.Range("I96").Value = 411.1
$ws.Range("J96").Value = 1115.5
$ws.Range("K96").Value = 1233.3
$ws.Range("L96").Value = 3346.5
$ws.Range("M96").Value = 139.6999999999998
$ws.Range("N96").Value = -6092.5
$ws.Range("H103").Value = 3820
$ws.Range("I103").Value = 2587.4443
$ws.Range("J103").Value = 4828.4546
$ws.Range("K103").Value = 7762.3329
$ws.Range("L103").Value = 14485.3638
$ws.Range("M103").Value = -7176.3329
$ws.Range("N103").Value = -15657.3638
$ws.Range("H129").Value = 2118.2
$ws.Range("J129").Value = 2556.5715
$ws.Range("L129").Value = 7669.7145
$ws.Range("N129").Value = -17669.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 536.9474
$ws.Range("I97").Value = 536.9474
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 536.9474
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -40.94740000000002
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 3149.9375
$ws.Range("I132").Value = 1310
$ws.Range("K132").Value = 3930
$ws.Range("M132").Value = -1400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 2531.625
$ws.Range("I37").Value = 2531.625
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2531.625
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2394.625
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3019.7144
$ws.Range("I22").Value = 379.66666
$ws.Range("K22").Value = 379.66666
$ws.Range("M22").Value = -29.66665999999998
$ws.Range("H68").Value = 71220
$ws.Range("J68").Value = 71220
$ws.Range("L68").Value = 71220
$ws.Range("N68").Value = -72718
$ws.Range("H71").Value = 71220
$ws.Range("J71").Value = 71220
$ws.Range("L71").Value = 213660
$ws.Range("N71").Value = -221148
$ws.Range("H94").Value = 4357.3
$ws.Range("I94").Value = 2012
$ws.Range("K94").Value = 2012
$ws.Range("M94").Value = -1561
$ws.Range("H99").Value = 3790.4783
$ws.Range("I99").Value = 3465.4119
$ws.Range("J99").Value = 4711.5
$ws.Range("K99").Value = 3465.4119
$ws.Range("L99").Value = 4711.5
$ws.Range("M99").Value = -1967.4119
$ws.Range("N99").Value = -7707.5
$ws.Range("H126").Value = 3790.4783
$ws.Range("I126").Value = 3465.4119
$ws.Range("J126").Value = 4711.5
$ws.Range("K126").Value = 10396.2357
$ws.Range("L126").Value = 14134.5
$ws.Range("M126").Value = -7926.235700000001
$ws.Range("N126").Value = -19074.5
$ws.Range("H132").Value = 3506.5833
$ws.Range("I132").Value = 3150.625
$ws.Range("J132").Value = 4218.5
$ws.Range("K132").Value = 9451.875
$ws.Range("L132").Value = 12655.5
$ws.Range("M132").Value = -6921.875
$ws.Range("N132").Value = -17715.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 29.565218
$ws.Range("I2").Value = 28.454546
$ws.Range("J2").Value = 30.583334
$ws.Range("K2").Value = 170.727276
$ws.Range("L2").Value = 183.500004
$ws.Range("M2").Value = -57.72727600000002
$ws.Range("N2").Value = -409.500004
$ws.Range("H34").Value = 3123.6667
$ws.Range("I34").Value = 144
$ws.Range("J34").Value = 3496.125
$ws.Range("K34").Value = 432
$ws.Range("L34").Value = 10488.375
$ws.Range("N34").Value = -10656.375
$ws.Range("M34").Value = -348
$ws.Range("H39").Value = 6805.6113
$ws.Range("J39").Value = 8606
$ws.Range("L39").Value = 25818
$ws.Range("N39").Value = -26406
$ws.Range("H121").Value = 220.8
$ws.Range("I121").Value = 220.8
$ws.Range("K121").Value = 662.4000000000001
$ws.Range("M121").Value = 647.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1793.0333
$ws.Range("I102").Value = 1436.7407
$ws.Range("K102").Value = 1436.7407
$ws.Range("M102").Value = 185.2592999999999
$ws.Range("H122").Value = 337576.47
$ws.Range("I122").Value = 388531.53
$ws.Range("K122").Value = 1165594.59
$ws.Range("M122").Value = -1163144.59

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 790
$ws.Range("I31").Value = 530
$ws.Range("J31").Value = 1050
$ws.Range("K31").Value = 530
$ws.Range("L31").Value = 1050
$ws.Range("M31").Value = -282
$ws.Range("N31").Value = -1546
$ws.Range("H40").Value = 6032.65
$ws.Range("I40").Value = 4725.8
$ws.Range("K40").Value = 4725.8
$ws.Range("M40").Value = -4589.8
$ws.Range("H93").Value = 2304.2
$ws.Range("I93").Value = 2226.889
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2226.889
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -978.8890000000001
$ws.Range("N93").Value = -5496
$ws.Range("H132").Value = 6198.8
$ws.Range("I132").Value = 5247.25
$ws.Range("J132").Value = 10005
$ws.Range("K132").Value = 15741.75
$ws.Range("L132").Value = 30015
$ws.Range("M132").Value = -13211.75
$ws.Range("N132").Value = -35075
$ws.Range("H136").Value = 4924
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1618
$ws.Range("J96").Value = 1156.3334
$ws.Range("L96").Value = 1156.3334
$ws.Range("N96").Value = -3902.3334
